$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/16/2023  Through  1/22/2023"

# --- Data table updates (rows 15-27) ---
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 2
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = -90
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 19
$ws.Range("H16").Value = -42.105263157894
$ws.Range("I16").Value = 11
$ws.Range("J16").Value = 16
$ws.Range("K16").Value = -31.25
$ws.Range("L16").Value = -8.333333333333
$ws.Range("M16").Value = -38.888888888888
$ws.Range("N16").Value = -92.361111111111
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 30
$ws.Range("G17").Value = 31
$ws.Range("H17").Value = -3.225806451612
$ws.Range("I17").Value = 23
$ws.Range("J17").Value = 22
$ws.Range("K17").Value = 4.545454545454
$ws.Range("L17").Value = 109.090909090909
$ws.Range("M17").Value = 27.777777777777
$ws.Range("N17").Value = -50
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 60
$ws.Range("I18").Value = 11
$ws.Range("J18").Value = 8
$ws.Range("K18").Value = 37.5
$ws.Range("L18").Value = 10
$ws.Range("M18").Value = -21.428571428571
$ws.Range("N18").Value = -95.045045045045
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -9.090909090909
$ws.Range("F19").Value = 52
$ws.Range("G19").Value = 46
$ws.Range("H19").Value = 13.043478260869
$ws.Range("I19").Value = 41
$ws.Range("J19").Value = 37
$ws.Range("K19").Value = 10.810810810810
$ws.Range("L19").Value = 86.363636363636
$ws.Range("M19").Value = 2.5
$ws.Range("N19").Value = -8.888888888888
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -75
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 8
$ws.Range("J20").Value = 7
$ws.Range("K20").Value = 14.285714285714
$ws.Range("L20").Value = -20
$ws.Range("M20").Value = -42.857142857142
$ws.Range("N20").Value = -95.428571428571
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 33
$ws.Range("E21").Value = -27.272727272727
$ws.Range("F21").Value = 121
$ws.Range("G21").Value = 118
$ws.Range("H21").Value = 2.542372881355
$ws.Range("I21").Value = 94
$ws.Range("J21").Value = 92
$ws.Range("K21").Value = 2.173913043478
$ws.Range("L21").Value = 40.298507462686
$ws.Range("M21").Value = -10.476190476190
$ws.Range("N21").Value = -85.381026438569
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = "0"
$ws.Range("C22").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("C24").Value = 40
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = 53.846153846153
$ws.Range("F24").Value = 117
$ws.Range("G24").Value = 91
$ws.Range("H24").Value = 28.571428571428
$ws.Range("I24").Value = 83
$ws.Range("J24").Value = 70
$ws.Range("K24").Value = 18.571428571428
$ws.Range("L24").Value = 6.410256410256
$ws.Range("M24").Value = 6.410256410256
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 42
$ws.Range("G25").Value = 38
$ws.Range("H25").Value = 10.526315789473
$ws.Range("I25").Value = 30
$ws.Range("J25").Value = 27
$ws.Range("K25").Value = 11.111111111111
$ws.Range("L25").Value = 7.142857142857
$ws.Range("M25").Value = -28.571428571428
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 1
$ws.Range("E26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 33.333333333333
$ws.Range("I26").Value = 2
$ws.Range("J26").Value = 2
$ws.Range("L26").Value = -50
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 2
$ws.Range("J27").Value = 3
$ws.Range("K27").Value = -33.333333333333
$ws.Range("L27").Value = -33.333333333333

$excel.CutCopyMode = 0
